$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 6438
$ws.Range("F4").Value = 1031
$ws.Range("F5").Value = 630
$ws.Range("F6").Value = 1417
$ws.Range("F7").Value = 3172
$ws.Range("F8").Value = 5
$ws.Range("F9").Value = 519
$ws.Range("F10").Value = 2064
$ws.Range("F13").Value = 216
$ws.Range("F14").Value = 104
$ws.Range("F15").Value = 221
$ws.Range("F16").Value = 1026
$ws.Range("F17").Value = 389
$ws.Range("F19").Value = 144
$ws.Range("F20").Value = 3934
$ws.Range("F22").Value = 3105
$ws.Range("F23").Value = 301
$ws.Range("F24").Value = 72
$ws.Range("F25").Value = 2670
$ws.Range("F26").Value = 2671
$ws.Range("F27").Value = 4499
$ws.Range("F29").Value = 948
$ws.Range("F30").Value = 501
$ws.Range("F31").Value = 2981
$ws.Range("F32").Value = 268
$ws.Range("F33").Value = 36
$ws.Range("F34").Value = 93
$ws.Range("F35").Value = 58
$ws.Range("F37").Value = 1082
$ws.Range("F38").Value = 1332
$ws.Range("F39").Value = 93
$ws.Range("F40").Value = 1188
$ws.Range("F41").Value = 769
$ws.Range("F42").Value = 4
$ws.Range("F43").Value = 704
$ws.Range("F44").Value = 464
$ws.Range("F46").Value = 168
$ws.Range("F47").Value = 20
$ws.Range("F48").Value = 66
$ws.Range("F49").Value = 337
$ws.Range("F50").Value = 3649

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F10").Value = 942
$ws.Range("F19").Value = 5
$ws.Range("F25").Value = 24

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1185

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 6438
$ws.Range("F5").Value = 630
$ws.Range("F6").Value = 1417
$ws.Range("F7").Value = 3172
$ws.Range("F8").Value = 5
$ws.Range("F9").Value = 519
$ws.Range("F11").Value = 2064
$ws.Range("F14").Value = 216
$ws.Range("F15").Value = 942
$ws.Range("F17").Value = 104
$ws.Range("F18").Value = 221
$ws.Range("F19").Value = 1026
$ws.Range("F21").Value = 389
$ws.Range("F22").Value = 144
$ws.Range("F23").Value = 3934
$ws.Range("F27").Value = 3105
$ws.Range("F28").Value = 2671
$ws.Range("F29").Value = 2671
$ws.Range("F30").Value = 4499
$ws.Range("F31").Value = 948
$ws.Range("F32").Value = 2981
$ws.Range("F33").Value = 268
$ws.Range("F35").Value = 1082
$ws.Range("F36").Value = 1332
$ws.Range("F37").Value = 93
$ws.Range("F38").Value = 1188
$ws.Range("F39").Value = 769
$ws.Range("F41").Value = 464
$ws.Range("F45").Value = 24
$ws.Range("F46").Value = 169
$ws.Range("F47").Value = 20
$ws.Range("F48").Value = 66
$ws.Range("F49").Value = 337
$ws.Range("F50").Value = 3649
